# 9.1.2 — add a new "2022" column (S) to the transport stats table,
# copying the formatting from the existing "2021" column (R) and
# filling in the new year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats — used so the new column inherits the same visual
# formatting (number format / font / borders / alignment) as column R,
# without creating redundant style entries.
$xlPasteFormats = -4122

function Copy-ColumnFormat($row) {
    $ws.Range("R$row").Copy() | Out-Null
    $ws.Range("S$row").PasteSpecial($xlPasteFormats) | Out-Null
}

# Row 3 — header year
Copy-ColumnFormat 3
$ws.Range("S3").Value = 2022

# Row 4 — Перевозки пассажиров всеми видами транспорта
Copy-ColumnFormat 4
$ws.Range("S4").Value = 10444.200000000001

# Row 5 — Железнодорожный
Copy-ColumnFormat 5
$ws.Range("S5").Value = 21.7

# Row 6 — Автомобильный
Copy-ColumnFormat 6
$ws.Range("S6").Value = 7361.6

# Row 7 — Водный транспорт
Copy-ColumnFormat 7
$ws.Range("S7").Value = 143.1

# Row 8 — Воздушный транспорт
Copy-ColumnFormat 8
$ws.Range("S8").Value = 844.2

# Row 9 — Трубопроводный транспорт (no data yet for 2022 — leave blank)
Copy-ColumnFormat 9

# Row 10 — Перевозки грузов всеми видами транспорта
Copy-ColumnFormat 10
$ws.Range("S10").Value = "2 756,0"

# Row 11 — Железнодорожный
Copy-ColumnFormat 11
$ws.Range("S11").Value = "1 013,8"

# Row 12 — Автомобильный
Copy-ColumnFormat 12
$ws.Range("S12").Value = "1 451,1"

# Row 13 — Водный транспорт
Copy-ColumnFormat 13
$ws.Range("S13").Value = 273.39999999999998

# Row 14 — Воздушный транспорт
Copy-ColumnFormat 14
$ws.Range("S14").Value = "-"

# Row 15 — Трубопроводный транспорт
Copy-ColumnFormat 15
$ws.Range("S15").Value = 17.7

# Match the saved selection state recorded in the workbook (cursor
# parked one cell to the right of the new column).
$ws.Range("T3").Select() | Out-Null
